# Daily attendance processing - 2026-02-01 19:20:39
#
# The "Recorded By" column (column G) lists the people who recorded a
# session's attendance as a comma-separated string. For every row where
# that value is currently "System, dnasr281@gmail.com", flip the order
# of the two names so it reads "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$recordedByCol = 7   # column G

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
